# "Added other user for myself"
# Adds a new form-response row (row 71), duplicating the last existing
# response (row 70) except for the relationship-type answer (column H),
# widens column G slightly, enlarges the header row + wraps the H1 header,
# and bumps the window zoom / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row: taller row + wrap text on the H1 header cell
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 228
$ws.Range("H1").WrapText = $true

# ---------------------------------------------------------------------
# 2) Column G gets a slightly different (custom) width
# ---------------------------------------------------------------------
$ws.Columns("G:G").ColumnWidth = 17.66

# ---------------------------------------------------------------------
# 3) New response row (row 71) - same data as row 70 except column H
# ---------------------------------------------------------------------
$ws.Range("A71").Value = 45340.566548159724
$ws.Range("B71").Value = "sebysurfer@gmail.com"
$ws.Range("C71").Value = "Sebastián W"
$ws.Range("D71").Value = "sebysurfer2003"
$ws.Range("E71").Value = 6241579275
$ws.Range("F71").Value = "M"
$ws.Range("G71").Value = "F"
$ws.Range("H71").Value = "DCA"
$ws.Range("I71").Value = 2
$ws.Range("J71").Value = "ICE"
$ws.Range("K71").Value = "COM"
$ws.Range("L71").Value = "TYG"
$ws.Range("M71").Value = -2
$ws.Range("N71").Value = 2
$ws.Range("O71").Value = -1
$ws.Range("P71").Value = 2
$ws.Range("Q71").Value = -2
$ws.Range("R71").Value = 2
$ws.Range("S71").Value = -2
$ws.Range("T71").Value = 0
$ws.Range("U71").Value = 1
$ws.Range("V71").Value = 1
$ws.Range("W71").Value = 1
$ws.Range("X71").Value = 2
$ws.Range("Y71").Value = 1
$ws.Range("Z71").Value = 2
$ws.Range("AA71").Value = 1
$ws.Range("AB71").Value = 2
$ws.Range("AC71").Value = -1
$ws.Range("AD71").Value = 2
$ws.Range("AE71").Value = 0
$ws.Range("AF71").Value = 2
$ws.Range("AG71").Value = 2
$ws.Range("AH71").Value = 2
$ws.Range("AI71").Value = -2
$ws.Range("AJ71").Value = 2
$ws.Range("AK71").Value = 2
$ws.Range("AL71").Value = 0
$ws.Range("AM71").Value = -1
$ws.Range("AN71").Value = -1
$ws.Range("AO71").Value = 0
$ws.Range("AP71").Value = -2
$ws.Range("AQ71").Value = -1

# Carry over row 70's formatting (number formats / fonts) onto row 71 in
# one shot, same as Excel does when a user fills a new row below an
# existing one.
$ws.Range("A70:AQ70").Copy()
$ws.Range("A71:AQ71").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Window / selection state
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 150
$ws.Range("H73").Select()
